$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (Exhibitions) - update "want to go" counts (column F)
# ---------------------------------------------------------------------
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 316
$wsExhibit.Range("F3").Value = 61
$wsExhibit.Range("F5").Value = 4616
$wsExhibit.Range("F6").Value = 358
$wsExhibit.Range("F8").Value = 285
$wsExhibit.Range("F9").Value = 716
$wsExhibit.Range("F10").Value = 199

# ---------------------------------------------------------------------
# Sheet "演出" (Performances) - append new event as row 3
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("A3").Value = 2
$wsShow.Range("B3").Value = "'2024-10-26"
$wsShow.Range("B3").Style = "Normal"
$wsShow.Range("C3").Value = "合肥·《四月是你的谎言》—“公生”与“薰”的钢琴小提琴唯美经典音乐集"
$wsShow.Range("D3").Value = "徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院"
$wsShow.Range("E3").Value = "2024.10.26 19:30-10.26 21:00"
$wsShow.Range("F3").Value = 0
$wsShow.Range("G3").Value = 40
$wsShow.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=90322"
$wsShow.Range("I3").Value = "//i2.hdslb.com/bfs/openplatform/202408/BiVgXUKH1722824304648.jpeg"
$wsShow.Range("A2").Copy()
$wsShow.Range("A3").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Sheet "全部类型" (All types) - same F-column updates as 展览, plus
# append the same new event as row 12
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 316
$wsAll.Range("F3").Value = 61
$wsAll.Range("F5").Value = 4616
$wsAll.Range("F6").Value = 358
$wsAll.Range("F8").Value = 285
$wsAll.Range("F9").Value = 716
$wsAll.Range("F11").Value = 199

$wsAll.Range("A12").Value = 11
$wsAll.Range("B12").Value = "'2024-10-26"
$wsAll.Range("B12").Style = "Normal"
$wsAll.Range("C12").Value = "合肥·《四月是你的谎言》—“公生”与“薰”的钢琴小提琴唯美经典音乐集"
$wsAll.Range("D12").Value = "徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院"
$wsAll.Range("E12").Value = "2024.10.26 19:30-10.26 21:00"
$wsAll.Range("F12").Value = 0
$wsAll.Range("G12").Value = 40
$wsAll.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=90322"
$wsAll.Range("I12").Value = "//i2.hdslb.com/bfs/openplatform/202408/BiVgXUKH1722824304648.jpeg"
$wsAll.Range("A11").Copy()
$wsAll.Range("A12").PasteSpecial(-4122)
